$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H2").Value = 2.86
$ws.Range("N2").Value = 3.9
$ws.Range("O2").Value = 1.32
$ws.Range("AK2").Value = 30
$ws.Range("AM2").Value = 100
$ws.Range("AN2").Value = 24
$ws.Range("G3").Value = 2.8
$ws.Range("H3").Value = 2.66
$ws.Range("L3").Value = 1.36
$ws.Range("N3").Value = 4.3
$ws.Range("O3").Value = 1.27
$ws.Range("P3").Value = 2.14
$ws.Range("Q3").Value = 1.78
$ws.Range("R3").Value = 1.45
$ws.Range("S3").Value = 3.05
$ws.Range("T3").Value = 1.66
$ws.Range("U3").Value = 2.36
$ws.Range("W3").Value = 1.55
$ws.Range("X3").Value = 21
$ws.Range("AB3").Value = 14.5
$ws.Range("AF3").Value = 23
$ws.Range("AI3").Value = 44
$ws.Range("AM3").Value = 85
$ws.Range("F4").Value = 5.1
$ws.Range("I4").Value = 1.91
$ws.Range("J4").Value = 3.4
$ws.Range("K4").Value = 4
$ws.Range("N4").Value = 2.86
$ws.Range("P4").Value = 1.62
$ws.Range("R4").Value = 1.23
$ws.Range("T4").Value = 2.06
$ws.Range("V4").Value = 2.1
$ws.Range("AA4").Value = 25
$ws.Range("AB4").Value = 18
$ws.Range("AD4").Value = 11.5
$ws.Range("AG4").Value = 28
$ws.Range("AH4").Value = 26
$ws.Range("F5").Value = 2.7
$ws.Range("G5").Value = 3.15
$ws.Range("H5").Value = 3
$ws.Range("O5").Value = 1.66
$ws.Range("W5").Value = 1.46
$ws.Range("Y5").Value = 8.6
$ws.Range("F8").Value = 8.6
$ws.Range("J8").Value = 4.4
$ws.Range("R8").Value = 1.39
$ws.Range("V8").Value = 3.3
$ws.Range("H9").Value = 3.1
$ws.Range("J9").Value = 1.03
$ws.Range("Q9").Value = 1.02
$ws.Range("W9").Value = 1.05
$ws.Range("F10").Value = 1.04
$ws.Range("G10").Value = 1000
$ws.Range("H10").Value = 1.04
$ws.Range("I10").Value = 1000
$ws.Range("J10").Value = 1.04
$ws.Range("K10").Value = 1000
$ws.Range("N10").Value = 1.1
$ws.Range("Q10").Value = 1.01
$ws.Range("S10").Value = 1.05
$ws.Range("V10").Value = 1.01
$ws.Range("W10").Value = 1.01
$ws.Range("K11").Value = 4.5
$ws.Range("N11").Value = 3.3
$ws.Range("O11").Value = 1.36
$ws.Range("P11").Value = 1.83
$ws.Range("Q11").Value = 2.02
$ws.Range("S11").Value = 3.75
$ws.Range("AB11").Value = 8.800000000000001
$ws.Range("AM11").Value = 200
$ws.Range("P12").Value = 2.36
$ws.Range("R12").Value = 1.57
$ws.Range("AE12").Value = 980
$ws.Range("AM12").Value = 70
$ws.Range("F13").Value = 2.24
$ws.Range("I13").Value = 3.5
$ws.Range("K13").Value = 3.75
$ws.Range("N13").Value = 4.3
$ws.Range("O13").Value = 1.23
$ws.Range("P13").Value = 2.14
$ws.Range("Q13").Value = 1.7
$ws.Range("R13").Value = 1.47
$ws.Range("S13").Value = 2.76
$ws.Range("U13").Value = 2.38
$ws.Range("W13").Value = 1.71
$ws.Range("AC13").Value = 8.800000000000001
$ws.Range("F15").Value = 1.76
$ws.Range("G15").Value = 1.92
$ws.Range("H15").Value = 4.1
$ws.Range("I15").Value = 5.4
$ws.Range("K15").Value = 5.3
$ws.Range("Q15").Value = 1.6
$ws.Range("V15").Value = 1.23
$ws.Range("W15").Value = 2.08
$ws.Range("F16").Value = 2.72
$ws.Range("G16").Value = 2.96
$ws.Range("P16").Value = 1.81
$ws.Range("R16").Value = 1.33
$ws.Range("W16").Value = 1.51
$ws.Range("AK16").Value = 40
$ws.Range("G17").Value = 4.4
$ws.Range("H17").Value = 1.95
$ws.Range("J17").Value = 3.7
$ws.Range("L17").Value = 1.27
$ws.Range("N17").Value = 4.5
$ws.Range("O17").Value = 1.22
$ws.Range("P17").Value = 2.22
$ws.Range("Q17").Value = 1.66
$ws.Range("R17").Value = 1.48
$ws.Range("S17").Value = 2.62
$ws.Range("T17").Value = 1.61
$ws.Range("U17").Value = 2.3
$ws.Range("V17").Value = 1.84
$ws.Range("W17").Value = 1.29
$ws.Range("X17").Value = 25
$ws.Range("Y17").Value = 14.5
$ws.Range("AJ17").Value = 85
$ws.Range("AL17").Value = 55
$ws.Range("AM17").Value = 85
$ws.Range("AN17").Value = 40
$ws.Range("AO17").Value = 14
$ws.Range("G19").Value = 2.34
$ws.Range("H19").Value = 3
$ws.Range("I19").Value = 3.45
$ws.Range("N19").Value = 5.2
$ws.Range("R19").Value = 1.57
$ws.Range("S19").Value = 2.42
$ws.Range("W19").Value = 1.74
$ws.Range("H20").Value = 9.199999999999999
$ws.Range("L20").Value = 1.19
$ws.Range("AF20").Value = 17
$ws.Range("F21").Value = 2.16
$ws.Range("I21").Value = 3.35
$ws.Range("K21").Value = 4.4
$ws.Range("N21").Value = 5.3
$ws.Range("P21").Value = 2.46
$ws.Range("Q21").Value = 1.56
$ws.Range("R21").Value = 1.6
$ws.Range("S21").Value = 2.24
$ws.Range("T21").Value = 1.53
$ws.Range("U21").Value = 2.52
$ws.Range("AB21").Value = 18
$ws.Range("AC21").Value = 12
$ws.Range("F22").Value = 5.2
$ws.Range("G22").Value = 9.4
$ws.Range("H22").Value = 1.42
$ws.Range("K22").Value = 6.8
$ws.Range("P22").Value = 2.62
$ws.Range("R22").Value = 1.62
$ws.Range("P23").Value = 2.32
$ws.Range("F24").Value = 4
$ws.Range("H24").Value = 1.85
$ws.Range("H25").Value = 6.4
$ws.Range("P25").Value = 2.64
$ws.Range("Q25").Value = 1.56
$ws.Range("AA25").Value = 190
$ws.Range("G26").Value = 3.25
$ws.Range("H26").Value = 2.64
$ws.Range("K26").Value = 3.2
$ws.Range("T26").Value = 1.87
$ws.Range("I28").Value = 2.66
$ws.Range("V28").Value = 1.6
$ws.Range("F29").Value = 2.68
$ws.Range("H29").Value = 2.86
$ws.Range("K29").Value = 3.35
$ws.Range("N29").Value = 2.88
$ws.Range("R29").Value = 1.23
$ws.Range("T29").Value = 1.93
$ws.Range("U29").Value = 1.92
$ws.Range("F30").Value = 3.85
$ws.Range("I30").Value = 2.06
$ws.Range("J30").Value = 3.75
$ws.Range("K30").Value = 4.4
$ws.Range("N30").Value = 4.2
$ws.Range("O30").Value = 1.24
$ws.Range("P30").Value = 2.12
$ws.Range("Q30").Value = 1.71
$ws.Range("R30").Value = 1.44
$ws.Range("S30").Value = 2.8
$ws.Range("T30").Value = 1.67
$ws.Range("U30").Value = 2.2
$ws.Range("V30").Value = 1.95
$ws.Range("X30").Value = 20
$ws.Range("AA30").Value = 24
$ws.Range("AB30").Value = 18.5
$ws.Range("AK30").Value = 980
$ws.Range("AN30").Value = 50
$ws.Range("AO30").Value = 12
$ws.Range("G31").Value = 2.1
$ws.Range("H31").Value = 4.3
$ws.Range("Q31").Value = 2.48
$ws.Range("S31").Value = 4.8
$ws.Range("W31").Value = 1.92
$ws.Range("F32").Value = 5.4
$ws.Range("G32").Value = 6.2
$ws.Range("H32").Value = 1.68
$ws.Range("J32").Value = 3.85
$ws.Range("K32").Value = 4.3
$ws.Range("N32").Value = 3.85
$ws.Range("O32").Value = 1.29
$ws.Range("P32").Value = 1.99
$ws.Range("Q32").Value = 1.84
$ws.Range("R32").Value = 1.38
$ws.Range("S32").Value = 3.1
$ws.Range("T32").Value = 1.83
$ws.Range("U32").Value = 2.04
$ws.Range("W32").Value = 1.19
$ws.Range("X32").Value = 19.5
$ws.Range("AI32").Value = 40
$ws.Range("AK32").Value = 80
$ws.Range("AL32").Value = 75
$ws.Range("AM32").Value = 120
$ws.Range("AN32").Value = 110
$ws.Range("AO32").Value = 10.5
